$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph "... ser-lhe-á mostrada as diversas informações acerca do
#    cliente/contador. " -> drop the trailing period + trailing space run,
#    then append the new parenthetical detail as a sequence of new runs.
# ---------------------------------------------------------------------------

$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*mostrada as diversas informa*cliente/contador*") {
        $targetIdx = $i
    }
}

# Remove the trailing " " run (whole-run delete keeps neighbouring runs intact)
$targetPara = $d.Paragraphs.Item($targetIdx)
$pEnd = $targetPara.Range.End
$spaceRun = $d.Range($pEnd - 2, $pEnd - 1)
$spaceRun.Delete()

# Remove the whole " diversas informações acerca do cliente/contador." run
$oldChunk = " diversas informações acerca do cliente/contador."
$targetPara = $d.Paragraphs.Item($targetIdx)
$pEnd2 = $targetPara.Range.End
$bigRun = $d.Range($pEnd2 - 1 - $oldChunk.Length, $pEnd2 - 1)
$bigRun.Delete()

# Re-append the replacement text as a series of distinct runs
$newChunks = @(
    " diversas informações acerca do cliente/contador",
    " (",
    "como nº de co",
    "ntador,",
    " morada, cliente, ",
    "data início de contrato e data fim de contrato",
    ", ocorrências",
    ")."
)
foreach ($chunk in $newChunks) {
    $targetPara = $d.Paragraphs.Item($targetIdx)
    $e = $targetPara.Range.End
    $insertionPoint = $d.Range($e - 1, $e - 1)
    $insertionPoint.InsertAfter($chunk)
}

# ---------------------------------------------------------------------------
# 2) The following (empty) paragraph becomes underlined-by-default for any
#    future text -- i.e. its paragraph mark gets <w:rPr><w:u val="single"/>.
# ---------------------------------------------------------------------------

$emptyIdx = $targetIdx + 1

$emptyPara = $d.Paragraphs.Item($emptyIdx)
$ip = $d.Range($emptyPara.Range.Start, $emptyPara.Range.Start)
$ip.InsertAfter("X")

$emptyPara = $d.Paragraphs.Item($emptyIdx)
$emptyPara.Range.Underline = 1

$emptyPara = $d.Paragraphs.Item($emptyIdx)
$tempCharRange = $d.Range($emptyPara.Range.Start, $emptyPara.Range.Start + 1)
$tempCharRange.Delete()

# ---------------------------------------------------------------------------
# 3) Drop the final "User stories-" paragraph entirely.
# ---------------------------------------------------------------------------

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Delete()
